$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells being edited keep plain-text formatting so Excel
# does not reinterpret dotted-number-looking strings (e.g. "525.43",
# "341.10", "0.990", "0.0665") as numeric values and mangle them via
# floating point / scientific notation round-tripping.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.752.24'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.086.59'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '525.43'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.74'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.70%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.085.59'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.442'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.14'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.79%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.391'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.615.56'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.42'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -6.82%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '57.790.86'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.084.33'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.08'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.65'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.17%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '341.10'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.31'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.90%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.60%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0913'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.37'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.22'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.89'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.71'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.99%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.14'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.92'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.15%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0665'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.27%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +7.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.97'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.683'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.126.72'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.89'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.276.29'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.990'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.63%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.49'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.98%  '
